# Weekly/daily price update for "Hortaliza, Vega Monumental Concepción - Tomate".
# Two new price records (Primera / Segunda grades) are published for 2023-10-13
# (Excel serial date 45212). They are inserted as new rows 817-818, pushing all
# existing records below down by two rows (817->819 ... 900->902), which is
# exactly what the diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 817, shifting rows 817:900 down
# to become rows 819:902 and extending the sheet dimension to A1:R902.
$ws.Range("A817:R818").Insert()

# ---- New row 817: Tomate, Larga vida, Primera --------------------------------
$ws.Cells.Item(817, 1).Value  = 11
$ws.Cells.Item(817, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(817, 3).Value  = "Bíobío"
$ws.Cells.Item(817, 4).Value  = 45212
$ws.Cells.Item(817, 5).Value  = 8
$ws.Cells.Item(817, 6).Value  = 100112020
$ws.Cells.Item(817, 7).Value  = "Tomate"
$ws.Cells.Item(817, 8).Value  = "Larga vida"
$ws.Cells.Item(817, 9).Value  = "Primera"
$ws.Cells.Item(817, 10).Value = 300
$ws.Cells.Item(817, 11).Value = 12000
$ws.Cells.Item(817, 12).Value = 12000
$ws.Cells.Item(817, 13).Value = 12000
$ws.Cells.Item(817, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(817, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(817, 16).Value = 667
$ws.Cells.Item(817, 17).Value = 18
$ws.Cells.Item(817, 18).Value = "Hortaliza"

# ---- New row 818: Tomate, Larga vida, Segunda --------------------------------
$ws.Cells.Item(818, 1).Value  = 11
$ws.Cells.Item(818, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(818, 3).Value  = "Bíobío"
$ws.Cells.Item(818, 4).Value  = 45212
$ws.Cells.Item(818, 5).Value  = 8
$ws.Cells.Item(818, 6).Value  = 100112020
$ws.Cells.Item(818, 7).Value  = "Tomate"
$ws.Cells.Item(818, 8).Value  = "Larga vida"
$ws.Cells.Item(818, 9).Value  = "Segunda"
$ws.Cells.Item(818, 10).Value = 400
$ws.Cells.Item(818, 11).Value = 10000
$ws.Cells.Item(818, 12).Value = 10000
$ws.Cells.Item(818, 13).Value = 10000
$ws.Cells.Item(818, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(818, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(818, 16).Value = 556
$ws.Cells.Item(818, 17).Value = 18
$ws.Cells.Item(818, 18).Value = "Hortaliza"
